$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet
$ws.Name = "Test"

# 2. Fix the UserName for TC003: amberqa@mailinator.com -> amberroad@mailinator.com
$ws.Range("C4").Value = "amberroad@mailinator.com"

# 3. Add the new "ExpectedResult" / "Status" columns (header row first, then data
#    in the same left-to-right-ish order the original author typed them so that
#    the shared-string table indices line up: E1, F1, E2, E4, E3, E5)
$ws.Range("E1").Value = "ExpectedResult"
$ws.Range("F1").Value = "Status"
$ws.Range("E2").Value = "amitqaMain - Admin"
$ws.Range("E4").Value = "AmberRoad - Sub Account Admin"
$ws.Range("E3").Value = "ParcelShipment - Sub Account Admin"
$ws.Range("E5").Value = "Mainamitba - Admin"

# 4. Apply a thin box border around the whole used range. The original A:D
#    columns (text-formatted cells) keep their text number format, while the
#    new E:F columns get the border without forcing a text format.
$ws.Range("A1:D5").Borders.LineStyle = 1
$ws.Range("E1:F1").NumberFormat = "@"
$ws.Range("E1:F1").Borders.LineStyle = 1
$ws.Range("E2:F5").Borders.LineStyle = 1

# 5. Size the new E column to fit its content (best-fit width, ~34.57 chars).
$ws.Columns.Item(5).ColumnWidth = 33.6

# 6. Match the saved selection/active cell.
$ws.Range("F2:F5").Select() | Out-Null
